$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2801.3
$ws.Range("I62").Value = 3002.875
$ws.Range("J62").Value = 1995
$ws.Range("K62").Value = 3002.875
$ws.Range("L62").Value = 1995
$ws.Range("M62").Value = -2378.875
$ws.Range("N62").Value = -3243

$ws.Range("H65").Value = 2801.3
$ws.Range("I65").Value = 3002.875
$ws.Range("J65").Value = 1995
$ws.Range("K65").Value = 15014.375
$ws.Range("L65").Value = 9975
$ws.Range("M65").Value = -11894.375
$ws.Range("N65").Value = -16215

$ws.Range("H96").Value = 7999.6665
$ws.Range("I96").Value = 7499.5
$ws.Range("J96").Value = 9000
$ws.Range("K96").Value = 22498.5
$ws.Range("L96").Value = 27000
$ws.Range("M96").Value = -21125.5
$ws.Range("N96").Value = -29746

$ws.Range("H116").Value = 3752.5
$ws.Range("I116").Value = 3005
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 3005
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = 437
$ws.Range("N116").Value = -11384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 402.2
$ws.Range("I4").Value = 252.75
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 252.75
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -136.75
$ws.Range("N4").Value = -1232

$ws.Range("H5").Value = 63.666668
$ws.Range("I5").Value = 48
$ws.Range("J5").Value = 95
$ws.Range("K5").Value = 48
$ws.Range("L5").Value = 95
$ws.Range("M5").Value = 64
$ws.Range("N5").Value = -319

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 63.666668
$ws.Range("I4").Value = 48
$ws.Range("J4").Value = 95
$ws.Range("K4").Value = 48
$ws.Range("L4").Value = 95
$ws.Range("M4").Value = 67
$ws.Range("N4").Value = -325

$ws.Range("H20").Value = 3348.5715
$ws.Range("I20").Value = 2036
$ws.Range("J20").Value = 4333
$ws.Range("K20").Value = 2036
$ws.Range("L20").Value = 4333
$ws.Range("M20").Value = -1789
$ws.Range("N20").Value = -4827

$ws.Range("H22").Value = 5999
$ws.Range("I22").Value = 5999
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5999
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -5826

$ws.Range("H76").Value = 14389.667
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 14389.667
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 14389.667
$ws.Range("N76").Value = -15019.667

$ws.Range("H79").Value = 14389.667
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 14389.667
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 14389.667
$ws.Range("N79").Value = -16573.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

$ws.Range("H23").Value = 4009
$ws.Range("I23").Value = 4009
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 4009
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -3769
$ws.Range("N23").ClearContents()

$ws.Range("H27").Value = 4009
$ws.Range("I27").Value = 4009
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4009
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3817
$ws.Range("N27").ClearContents()

$ws.Range("H43").Value = 19492.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19492.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19492.5
$ws.Range("N43").Value = -19860.5

$ws.Range("H93").Value = 30189.5
$ws.Range("I93").Value = 28127.857
$ws.Range("J93").Value = 35000
$ws.Range("K93").Value = 28127.857
$ws.Range("L93").Value = 35000
$ws.Range("M93").Value = -26255.857
$ws.Range("N93").Value = -38744

$ws.Range("H101").Value = 19492.5
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 19492.5
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 19492.5
$ws.Range("N101").Value = -25982.5

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 8824.25
$ws.Range("I122").Value = 8824.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 26472.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -24022.75
$ws.Range("N122").ClearContents()

$ws.Range("H134").Value = 8999
$ws.Range("I134").Value = 8999
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 26997
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -24462
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 33.153847
$ws.Range("I12").Value = 45.75
$ws.Range("J12").Value = 27.555555
$ws.Range("K12").Value = 137.25
$ws.Range("L12").Value = 82.666665
$ws.Range("M12").Value = 35.75
$ws.Range("N12").Value = -428.666665

$ws.Range("H131").Value = 3390.2727
$ws.Range("I131").Value = 6432.3335
$ws.Range("J131").Value = 2249.5
$ws.Range("K131").Value = 19297.0005
$ws.Range("L131").Value = 6748.5
$ws.Range("M131").Value = -14257.0005
$ws.Range("N131").Value = -16828.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 17366666
$ws.Range("I11").Value = 47500000
$ws.Range("J11").Value = 2300000
$ws.Range("K11").Value = 47500000
$ws.Range("L11").Value = 2300000
$ws.Range("M11").Value = -47499861
$ws.Range("N11").Value = -2300278

$ws.Range("H107").Value = 2970.2
$ws.Range("I107").Value = 3249
$ws.Range("J107").Value = 1855
$ws.Range("K107").Value = 3249
$ws.Range("L107").Value = 1855
$ws.Range("M107").Value = -1329
$ws.Range("N107").Value = -5695

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4856.6
$ws.Range("I22").Value = 4095
$ws.Range("J22").Value = 5999
$ws.Range("K22").Value = 4095
$ws.Range("L22").Value = 5999
$ws.Range("M22").Value = -3800
$ws.Range("N22").Value = -6589

$ws.Range("H27").Value = 4856.6
$ws.Range("I27").Value = 4095
$ws.Range("J27").Value = 5999
$ws.Range("K27").Value = 4095
$ws.Range("L27").Value = 5999
$ws.Range("M27").Value = -3988
$ws.Range("N27").Value = -6213

$ws.Range("H41").Value = 14937
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 14937
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 14937
$ws.Range("N41").Value = -15813

$ws.Range("H132").Value = 8591.25
$ws.Range("I132").Value = 8580.25
$ws.Range("J132").Value = 8624.25
$ws.Range("K132").Value = 25740.75
$ws.Range("L132").Value = 25872.75
$ws.Range("M132").Value = -23210.75
$ws.Range("N132").Value = -30932.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H138").Value = 60000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 60000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280
